# Apply updates to the 2324_team_totals workbook:
#  - Update several total_pim (column F) values and one total_goals (B24) value
#    with refreshed query results.
#  - Remove the now-unused duplicate "total_goals"/"total_pim" columns (K, L)
#    that were left over from an earlier comparison query.
#  - Clear the stray "s" note in Q14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated total_pim (column F) values per refreshed query results.
$ws.Range("F2").Value = 1064
$ws.Range("F3").Value = 805
$ws.Range("F4").Value = 754
$ws.Range("F5").Value = 777
$ws.Range("F6").Value = 669
$ws.Range("F7").Value = 607
$ws.Range("F8").Value = 716
$ws.Range("F11").Value = 540
$ws.Range("F12").Value = 686
$ws.Range("F13").Value = 762
$ws.Range("F14").Value = 1104
$ws.Range("F15").Value = 693
$ws.Range("F16").Value = 910
$ws.Range("F17").Value = 742
$ws.Range("F19").Value = 731
$ws.Range("F20").Value = 580
$ws.Range("F21").Value = 729
$ws.Range("F22").Value = 808
$ws.Range("F24").Value = 562
$ws.Range("F25").Value = 627
$ws.Range("F26").Value = 745
$ws.Range("F27").Value = 542
$ws.Range("F30").Value = 760

# Updated total_goals (column B) value.
$ws.Range("B24").Value = 245

# Clear the leftover note in Q14.
$ws.Range("Q14").ClearContents()

# Remove the leftover duplicate comparison columns K and L entirely.
$ws.Range("K1:L32").EntireColumn.Delete()

$ws.Columns.AutoFit()

# Move the selection to reflect where the analyst was working.
[void]$ws.Range("O9").Select()
